# Auto-generated edit script: update F-column (想去人数 / want-to-go counts)
# across all 4 worksheets to match the gh-pages snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5756
$ws.Range("F6").Value = 464
$ws.Range("F7").Value = 1046
$ws.Range("F8").Value = 3370
$ws.Range("F9").Value = 6588
$ws.Range("F10").Value = 199
$ws.Range("F11").Value = 1289
$ws.Range("F12").Value = 744
$ws.Range("F15").Value = 17
$ws.Range("F18").Value = 97
$ws.Range("F20").Value = 163
$ws.Range("F22").Value = 950
$ws.Range("F25").Value = 13
$ws.Range("F26").Value = 101
$ws.Range("F28").Value = 1149
$ws.Range("F30").Value = 30
$ws.Range("F33").Value = 280
$ws.Range("F35").Value = 267
$ws.Range("F36").Value = 1164
$ws.Range("F38").Value = 85

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F11").Value = 5
$ws.Range("F20").Value = 180
$ws.Range("F21").Value = 131
$ws.Range("F24").Value = 604
$ws.Range("F28").Value = 659
$ws.Range("F29").Value = 949
$ws.Range("F30").Value = 560
$ws.Range("F32").Value = 81
$ws.Range("F35").Value = 93
$ws.Range("F38").Value = 50

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 538
$ws.Range("F7").Value = 279
$ws.Range("F8").Value = 884

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F8").Value = 538
$ws.Range("F9").Value = 538
$ws.Range("F10").Value = 279
$ws.Range("F11").Value = 279
$ws.Range("F14").Value = 5756
$ws.Range("F15").Value = 464
$ws.Range("F16").Value = 1046
$ws.Range("F17").Value = 3370
$ws.Range("F19").Value = 6588
$ws.Range("F20").Value = 199
$ws.Range("F21").Value = 1289
$ws.Range("F24").Value = 744
$ws.Range("F26").Value = 884
$ws.Range("F27").Value = 180
$ws.Range("F28").Value = 17
$ws.Range("F30").Value = 97
$ws.Range("F31").Value = 163
$ws.Range("F32").Value = 950
$ws.Range("F33").Value = 604
$ws.Range("F35").Value = 13
$ws.Range("F36").Value = 101
$ws.Range("F37").Value = 1149
$ws.Range("F39").Value = 30
$ws.Range("F42").Value = 949
$ws.Range("F43").Value = 560
$ws.Range("F44").Value = 280
$ws.Range("F45").Value = 81
$ws.Range("F46").Value = 267
$ws.Range("F47").Value = 93
$ws.Range("F50").Value = 85
$ws.Range("F51").Value = 50
